$d = $word.ActiveDocument

# Remove the entire second "Pack Label" block: from the page break that
# precedes "2. AIS deikdjmewl232 Pack Label" through the final
# "Quantity: 100" paragraph of that block (inclusive). The page break that
# follows it (right before the section end) is left in place.
$secondBlockStart = $d.Paragraphs.Item(9).Range.Start
$secondBlockEnd = $d.Paragraphs.Item(17).Range.End
$d.Range($secondBlockStart, $secondBlockEnd).Delete()

# Apply the text edits for the remaining (first) label block.
$d.Content.Find.Execute("1. AIS Adkdm232 Pack Label", $true, $false, $false, $false, $false, $true, 1, $false, "1. asi ded Pack Label", 2)
$d.Content.Find.Execute("Part Number: Adkdm232", $true, $false, $false, $false, $false, $true, 1, $false, "Part Number: ded", 2)
$d.Content.Find.Execute("Description: djeidjoewd", $true, $false, $false, $false, $false, $true, 1, $false, "Description: ded", 2)
$d.Content.Find.Execute("PO Number: PO2333", $true, $false, $false, $false, $false, $true, 1, $false, "PO Number: ewdw", 2)
$d.Content.Find.Execute("Mfg. 04-2025", $true, $false, $false, $false, $false, $true, 1, $false, "Mfg. 08-2025", 2)
